# Updates the delivery demo sheet to match the new Batch record string
# representation (adds "| <id>" suffix, refreshed quantities/prices/dates),
# and introduces a second weekly delivery whose header row (date/code) had
# previously been collapsed into row 4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: existing delivery header + first batch line -------------------
$ws.Range("B2").Value = "'123"
$ws.Range("E2").Value = "(Aciclovir) 200 mg Tablet 5x: " + [char]8369 + "5.0 (2021-02-12 | 1)"

# --- Row 3: second batch line for the first delivery -----------------------
$ws.Range("E3").Value = "Mucosolve (Ambroxol) 15 mg/60 ml Syrup 1x: " + [char]8369 + "9.0 (2021-02-19 | 45)"

# --- Row 4: no longer a delivery header; becomes a plain batch line --------
$ws.Range("A4:D4").ClearContents()
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("E4").Value = "(Aciclovir) 200 mg Tablet 1x: " + [char]8369 + "2.0 (2021-02-15 | 4)"

# --- Row 5: new delivery header (date/code/from/to) + batch line ----------
$ws.Range("A2").Copy()
$ws.Range("A5").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("A5").Value = 44235
$ws.Range("B5").Value = "'43253"
$ws.Range("C5").Value = "My Company"
$ws.Range("D5").Value = "My Company"
$ws.Range("E5").Value = "(Allopurinol) 100 mg Tablet 10x: " + [char]8369 + "2.0 (2021-02-17 | 456)"

# --- Row 6: third batch line, now attached to the second delivery ---------
$ws.Range("E6").Value = "Cisflem (Carbo) 125 mg/60 ml Syrup 50x: " + [char]8369 + "9.0 (2021-02-24 | 5678)"
